# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on the zh-cn and de-de
# report sheets to reflect the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 3 corresponds to file
# 9ed7af9a-55b6-484b-ba97-0453c63e4a9a.cb489f38db85672d977d9ecc4c64c5938868e9f0.zh-cn.xlf
$wsZhCn.Range("E3").Value = "2016-03-18 16:41:05"
$wsZhCn.Range("H3").Value = "2016-03-18 16:41:21"

# de-de sheet: row 3 corresponds to file
# 9ed7af9a-55b6-484b-ba97-0453c63e4a9a.cb489f38db85672d977d9ecc4c64c5938868e9f0.de-de.xlf
$wsDeDe.Range("E3").Value = "2016-03-18 16:41:09"
$wsDeDe.Range("H3").Value = "2016-03-18 16:41:27"
